$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: insert the two brand-new leading columns (_airbyte_ab_id, _airbyte_emitted_at) ----
$ws.Range("A1:B1").EntireColumn.Insert()

# ---- Step 2: insert the two brand-new columns before the old last column (now T) ----
# (_airbyte_additional_properties, source_file_path) -- old "updated_at" column slides from T to V
$ws.Range("T1:U1").EntireColumn.Insert()

# ---- Step 3: give the four new header cells the same bold/bordered/centered style as the rest of row 1 ----
$ws.Range("C1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("T1").PasteSpecial(-4122)
$ws.Range("U1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Step 4: give column B (the new _airbyte_emitted_at data cells) the same datetime format as column V ----
$ws.Range("V2").Copy()
$ws.Range("B2:B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Step 5: header row text ----
$ws.Range("A1").Value = "_airbyte_ab_id"
$ws.Range("B1").Value = "_airbyte_emitted_at"
$ws.Range("C1").Value = "sys_code"
$ws.Range("D1").Value = "sys_init"
$ws.Range("E1").Value = "sys_proc_id"
$ws.Range("F1").Value = "sys_acc_numb"
$ws.Range("G1").Value = "sys_acc_sold"
$ws.Range("H1").Value = "sys_cen_code"
$ws.Range("I1").Value = "sys_dom_code"
$ws.Range("J1").Value = "sys_ica_code"
$ws.Range("K1").Value = "sys_reg_code"
$ws.Range("L1").Value = "sys_acqu_iden"
$ws.Range("M1").Value = "sys_corp_iden"
$ws.Range("N1").Value = "sys_corp_name"
$ws.Range("O1").Value = "sys_stat_code"
$ws.Range("P1").Value = "sys_acc_numb_mxp"
$ws.Range("Q1").Value = "sys_acqu_busi_iden"
$ws.Range("R1").Value = "sys_sett_curr_code"
$ws.Range("S1").Value = "sys_sett_inst_iden"
$ws.Range("T1").Value = "_airbyte_additional_properties"
$ws.Range("U1").Value = "source_file_path"
$ws.Range("V1").Value = "updated_at"

# ---- Step 6: data rows ----
# row 2
$ws.Range("A2").Value = "'20bcc3bf-d1ce-4fed-b546-532ee6bac5cb"
$ws.Range("B2").Value = 45510.3079196875
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "'VISA"
$ws.Range("E2").Value = 999
$ws.Range("F2").Value = "'Visa_acct"
$ws.Range("G2").Value = 1
# H2: already empty (inherited from old F2) -- left untouched
$ws.Range("I2").Value = "'W"
$ws.Range("J2").Value = "'9999"
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = "'999999"
$ws.Range("M2").Value = "'VISA"
$ws.Range("N2").Value = "'VISA INTERNATIONAL"
# O2: already empty (inherited from old M2) -- left untouched
$ws.Range("P2").Value = "'visa_mxp_acct"
$ws.Range("Q2").Value = "'999999"
$ws.Range("R2").Value = 28
# S2: already empty (inherited from old Q2) -- left untouched
$ws.Range("T2").Value = "'"
$ws.Range("U2").Value = "'s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/SYSTEM/2024_08_06_1722929004063_0.parquet"
$ws.Range("V2").Value = 45511.29528810951
# row 3
$ws.Range("A3").Value = "'00cf7371-44d0-4302-ae4c-256754a44366"
$ws.Range("B3").Value = 45510.3079196875
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "'MCS"
$ws.Range("E3").Value = 118
$ws.Range("F3").Value = "'cis_acct"
$ws.Range("G3").Value = 1
# H3: already empty (inherited from old F3) -- left untouched
$ws.Range("I3").Value = "'W"
$ws.Range("J3").Value = "'12888"
$ws.Range("K3").Value = 13
$ws.Range("L3").Value = "'538815"
$ws.Range("M3").Value = "'MCI"
$ws.Range("N3").Value = "'MASTERCARD INT"
# O3: already empty (inherited from old M3) -- left untouched
$ws.Range("P3").Value = "'cis_mxp_acct"
$ws.Range("Q3").Value = "'538815"
$ws.Range("R3").Value = 28
# S3: already empty (inherited from old Q3) -- left untouched
$ws.Range("T3").Value = "'"
$ws.Range("U3").Value = "'s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/SYSTEM/2024_08_06_1722929004063_0.parquet"
$ws.Range("V3").Value = 45511.29528810951
# row 4
$ws.Range("A4").Value = "'f7d43022-89b6-479f-8b17-e8617ce47673"
$ws.Range("B4").Value = 45510.3079196875
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "'AMEX"
$ws.Range("E4").Value = 209
$ws.Range("F4").Value = "'amex_acct"
$ws.Range("G4").Value = 1
# H4: already empty (inherited from old F4) -- left untouched
$ws.Range("I4").Value = "'W"
$ws.Range("J4").Value = "'4104"
# K4: already empty (inherited from old I4) -- left untouched
$ws.Range("L4").Value = "'371111"
$ws.Range("M4").Value = "'AMEX"
$ws.Range("N4").Value = "'American Express"
# O4: already empty (inherited from old M4) -- left untouched
$ws.Range("P4").Value = "'amex_mxp_acct"
$ws.Range("Q4").Value = "'371111"
$ws.Range("R4").Value = 28
# S4: already empty (inherited from old Q4) -- left untouched
$ws.Range("T4").Value = "'"
$ws.Range("U4").Value = "'s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/SYSTEM/2024_08_06_1722929004063_0.parquet"
$ws.Range("V4").Value = 45511.29528810951
# row 5
$ws.Range("A5").Value = "'352d866d-5e22-4c82-bb2d-a1303ed4ecb1"
$ws.Range("B5").Value = 45510.3079196875
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "'PULSE"
$ws.Range("E5").Value = 172
$ws.Range("F5").Value = "'01900935280022"
$ws.Range("G5").Value = 1000000
# H5: already empty (inherited from old F5) -- left untouched
$ws.Range("I5").Value = "'W"
$ws.Range("J5").Value = "'172672"
# K5: already empty (inherited from old I5) -- left untouched
$ws.Range("L5").Value = "'1726726723"
$ws.Range("M5").Value = "'PULSE"
$ws.Range("N5").Value = "'PULSE"
# O5: already empty (inherited from old M5) -- left untouched
$ws.Range("P5").Value = "'MXPPULSE999"
$ws.Range("Q5").Value = "'17267267"
$ws.Range("R5").Value = 28
# S5: already empty (inherited from old Q5) -- left untouched
$ws.Range("T5").Value = "'"
$ws.Range("U5").Value = "'s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/SYSTEM/2024_08_06_1722929004063_0.parquet"
$ws.Range("V5").Value = 45511.29528810951
# row 6
$ws.Range("A6").Value = "'066616b2-ef52-4679-bdb0-84ae4a53b92a"
$ws.Range("B6").Value = 45510.3079196875
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = "'UPI"
$ws.Range("E6").Value = 103
$ws.Range("F6").Value = "'123456"
$ws.Range("G6").Value = 1000000
# H6: already empty (inherited from old F6) -- left untouched
$ws.Range("I6").Value = "'W"
$ws.Range("J6").Value = "'10344"
$ws.Range("K6").Value = 20
$ws.Range("L6").Value = "'31870524"
$ws.Range("M6").Value = "'UPI"
$ws.Range("N6").Value = "'UNION PAY INTERNATIONAL"
# O6: already empty (inherited from old M6) -- left untouched
$ws.Range("P6").Value = "'MXPUPI999"
$ws.Range("Q6").Value = "'31870524"
$ws.Range("R6").Value = 28
# S6: already empty (inherited from old Q6) -- left untouched
$ws.Range("T6").Value = "'"
$ws.Range("U6").Value = "'s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/SYSTEM/2024_08_06_1722929004063_0.parquet"
$ws.Range("V6").Value = 45511.29528810951
# row 7
$ws.Range("A7").Value = "'c8ff3e94-190c-4769-94f7-b0febc2a217a"
$ws.Range("B7").Value = 45510.3079196875
$ws.Range("C7").Value = 9
$ws.Range("D7").Value = "'D"
# E7: already empty (inherited from old C7) -- left untouched
$ws.Range("F7").Value = "'22-33-44-55"
$ws.Range("G7").Value = 0
# H7: already empty (inherited from old F7) -- left untouched
$ws.Range("I7").Value = "'W"
# J7: already empty (inherited from old H7) -- left untouched
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = "'356999"
$ws.Range("M7").Value = "'DCI"
$ws.Range("N7").Value = "'DCI"
# O7: already empty (inherited from old M7) -- left untouched
$ws.Range("P7").Value = "'22-33-44-55"
$ws.Range("Q7").Value = "'356999"
$ws.Range("R7").Value = 28
# S7: already empty (inherited from old Q7) -- left untouched
$ws.Range("T7").Value = "'"
$ws.Range("U7").Value = "'s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/SYSTEM/2024_08_06_1722929004063_0.parquet"
$ws.Range("V7").Value = 45511.29528810951
# row 8
$ws.Range("A8").Value = "'ae6e7ac1-3b5b-4696-bc41-c985561eaf18"
$ws.Range("B8").Value = 45510.3079196875
$ws.Range("C8").Value = 5241
$ws.Range("D8").Value = "'NEPS"
# E8: already empty (inherited from old C8) -- left untouched
# F8: already empty (inherited from old D8) -- left untouched
# G8: already empty (inherited from old E8) -- left untouched
# H8: already empty (inherited from old F8) -- left untouched
# I8: already empty (inherited from old G8) -- left untouched
# J8: already empty (inherited from old H8) -- left untouched
# K8: already empty (inherited from old I8) -- left untouched
# L8: already empty (inherited from old J8) -- left untouched
$ws.Range("M8").Value = "'NEPS"
$ws.Range("N8").Value = "'NEPS SWITCH"
# O8: already empty (inherited from old M8) -- left untouched
# P8: already empty (inherited from old N8) -- left untouched
# Q8: already empty (inherited from old O8) -- left untouched
# R8: already empty (inherited from old P8) -- left untouched
# S8: already empty (inherited from old Q8) -- left untouched
$ws.Range("T8").Value = "'"
$ws.Range("U8").Value = "'s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/SYSTEM/2024_08_06_1722929004063_0.parquet"
$ws.Range("V8").Value = 45511.29528810951
